$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "Nothing special to mention"; New = "Design: Nothing special to mention" },
    @{ Old = "91xx Went well, some minor budget challenges"; New = "Design: 91xx Went well, some minor budget challenges" },
    @{ Old = "Internal communication ok. External communication with suppliers mostly ok."; New = "Design: Internal communication ok. External communication with suppliers mostly ok." },
    @{ Old = "Some things went to correct direction but regarding TK the opposite way."; New = "Design: Some things went to correct direction but regarding TK the opposite way." },
    @{ Old = "Mostly"; New = "Design: Mostly" }
)

foreach ($rep in $replacements) {
    $d.Content.Find.Execute($rep.Old, $true, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2) | Out-Null
}
